{"js": "// Section 3 (\"Appointment Creation\") final sentence currently reads:\n//   \"...Send() follows the Wuu and Bernstein algorithm.\"\n// It should be changed to:\n//   \"...Send(), insert(), and delete() follow the Wuu and Bernstein algorithm.\"\n//\n// Do the edit as two small, in-place text replacements that line up with the\n// existing run boundaries, so only the words that actually changed are\n// touched (everything else, e.g. the \"Wuu\" spell-check run, is left alone).\n\nconst body = context.document.body;\n\n// 1) \"Send(\" -> \"Send()\"\nconst openParen = body.search(\"Send(\", { matchCase: true, matchWholeWord: false });\nopenParen.load(\"items\");\nawait context.sync();\n\nif (openParen.items.length > 0) {\n  openParen.items[0].insertText(\"Send()\", \"Replace\");\n  await context.sync();\n}\n\n// 2) \") follows the \" -> \", insert(), and delete() follow the \"\nconst tail = body.search(\") follows the \", { matchCase: true, matchWholeWord: false });\ntail.load(\"items\");\nawait context.sync();\n\nif (tail.items.length > 0) {\n  tail.items[0].insertText(\", insert(), and delete() follow the \", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Section 3 (\"Appointment Creation\") final sentence currently reads:\n#   \"...Send() follows the Wuu and Bernstein algorithm.\"\n# It should be changed to:\n#   \"...Send(), insert(), and delete() follow the Wuu and Bernstein algorithm.\"\n#\n# Do the edit as two small, in-place Find/Replace operations that line up\n# with the existing run boundaries, so only the words that actually changed\n# are touched.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n# 1) \"Send(\" -> \"Send()\"\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Execute(\n    \"Send(\",      # FindText\n    $true,        # MatchCase\n    $false,       # MatchWholeWord\n    $false,       # MatchWildcards\n    $false,       # MatchSoundsLike\n    $false,       # MatchAllWordForms\n    $true,        # Forward\n    $wdFindContinue,  # Wrap\n    $false,       # Format\n    \"Send()\",     # ReplaceWith\n    $wdReplaceOne # Replace\n) | Out-Null\n\n# 2) \") follows the \" -> \", insert(), and delete() follow the \"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Execute(\n    \") follows the \",\n    $true,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    $wdFindContinue,\n    $false,\n    \", insert(), and delete() follow the \",\n    $wdReplaceOne\n) | Out-Null\n"}
